$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 96.59999999999999
$ws.Range("I6").Value = 96.59999999999999
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 289.8
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -177.8
$ws.Range("N6").ClearContents()
$ws.Range("H12").Value = 208.83333
$ws.Range("I12").Value = 220.6
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 220.6
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = -50.59999999999999
$ws.Range("N12").Value = -490
$ws.Range("H33").Value = 500471.9
$ws.Range("I33").Value = 833599.9399999999
$ws.Range("K33").Value = 833599.9399999999
$ws.Range("M33").Value = -833370.9399999999
$ws.Range("H53").Value = 3055.2778
$ws.Range("J53").Value = 147.57143
$ws.Range("L53").Value = 147.57143
$ws.Range("N53").Value = -1421.57143
$ws.Range("H100").Value = 3349.5386
$ws.Range("I100").Value = 3043.6
$ws.Range("K100").Value = 3043.6
$ws.Range("M100").Value = -2502.6
$ws.Range("H107").Value = 1576.0834
$ws.Range("I107").Value = 1401.1818
$ws.Range("K107").Value = 1401.1818
$ws.Range("M107").Value = 518.8181999999999
$ws.Range("H111").Value = 3165.5
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 3165.5
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 9496.5
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -15630.5
$ws.Range("H118").Value = 2558.2727
$ws.Range("I118").Value = 2471.5557
$ws.Range("J118").Value = 2948.5
$ws.Range("K118").Value = 7414.6671
$ws.Range("L118").Value = 8845.5
$ws.Range("M118").Value = -5757.6671
$ws.Range("N118").Value = -12159.5
$ws.Range("H129").Value = 2542.111
$ws.Range("J129").Value = 3241.6667
$ws.Range("L129").Value = 9725.000100000001
$ws.Range("N129").Value = -19725.0001
$ws.Range("H136").Value = 152800
$ws.Range("J136").Value = 152800
$ws.Range("L136").Value = 152800
$ws.Range("N136").Value = -163000
$ws.Range("H137").Value = 945684.75
$ws.Range("I137").Value = 1033.1666
$ws.Range("J137").Value = 1418010.6
$ws.Range("K137").Value = 3099.4998
$ws.Range("L137").Value = 4254031.800000001
$ws.Range("M137").Value = -549.4998000000001
$ws.Range("N137").Value = -4259131.800000001
$ws.Range("H138").Value = 2245.35
$ws.Range("I138").Value = 2046.8
$ws.Range("J138").Value = 2311.5334
$ws.Range("K138").Value = 6140.4
$ws.Range("L138").Value = 6934.600199999999
$ws.Range("M138").Value = -1000.4
$ws.Range("N138").Value = -17214.6002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1549.0834
$ws.Range("I45").Value = 1513.4286
$ws.Range("J45").Value = 1599
$ws.Range("K45").Value = 1513.4286
$ws.Range("L45").Value = 1599
$ws.Range("M45").Value = -1136.4286
$ws.Range("N45").Value = -2353
$ws.Range("H61").Value = 4559.727
$ws.Range("I61").Value = 6395
$ws.Range("K61").Value = 6395
$ws.Range("M61").Value = -6183
$ws.Range("H104").Value = 38404
$ws.Range("J104").Value = 38404
$ws.Range("L104").Value = 38404
$ws.Range("N104").Value = -45392
$ws.Range("H110").Value = 1213.5714
$ws.Range("I110").Value = 1080.75
$ws.Range("J110").Value = 2010.5
$ws.Range("K110").Value = 1080.75
$ws.Range("L110").Value = 2010.5
$ws.Range("M110").Value = 964.25
$ws.Range("N110").Value = -6100.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 8479.883
$ws.Range("I132").Value = 10881.909
$ws.Range("J132").Value = 4076.1667
$ws.Range("K132").Value = 32645.727
$ws.Range("L132").Value = 12228.5001
$ws.Range("M132").Value = -30115.727
$ws.Range("N132").Value = -17288.5001
$ws.Range("H136").Value = 4559.727
$ws.Range("I136").Value = 6395
$ws.Range("K136").Value = 19185
$ws.Range("M136").Value = -16635

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2474.9375
$ws.Range("I86").Value = 2224.4614
$ws.Range("J86").Value = 3560.3333
$ws.Range("K86").Value = 2224.4614
$ws.Range("L86").Value = 3560.3333
$ws.Range("M86").Value = -1101.4614
$ws.Range("N86").Value = -5806.3333
$ws.Range("H89").Value = 2474.9375
$ws.Range("I89").Value = 2224.4614
$ws.Range("J89").Value = 3560.3333
$ws.Range("K89").Value = 11122.307
$ws.Range("L89").Value = 17801.6665
$ws.Range("M89").Value = -5506.307000000001
$ws.Range("N89").Value = -29033.6665
$ws.Range("H94").Value = 2243.5144
$ws.Range("I94").Value = 2194.724
$ws.Range("K94").Value = 2194.724
$ws.Range("M94").Value = -1743.724

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1184.625
$ws.Range("I16").Value = 1003
$ws.Range("J16").Value = 1293.6
$ws.Range("K16").Value = 1003
$ws.Range("L16").Value = 1293.6
$ws.Range("M16").Value = -716
$ws.Range("N16").Value = -1867.6
$ws.Range("H58").Value = 3214.9
$ws.Range("J58").Value = 3101.3333
$ws.Range("L58").Value = 3101.3333
$ws.Range("N58").Value = -3507.3333
$ws.Range("H86").Value = 90913690
$ws.Range("I86").Value = 125004330
$ws.Range("K86").Value = 125004330
$ws.Range("M86").Value = -125003207
$ws.Range("H89").Value = 90913690
$ws.Range("I89").Value = 125004330
$ws.Range("K89").Value = 625021650
$ws.Range("M89").Value = -625016034
$ws.Range("H99").Value = 2525.25
$ws.Range("I99").Value = 2431.077
$ws.Range("J99").Value = 2933.3333
$ws.Range("K99").Value = 2431.077
$ws.Range("L99").Value = 2933.3333
$ws.Range("M99").Value = -933.0770000000002
$ws.Range("N99").Value = -5929.3333
$ws.Range("H113").Value = 1184.625
$ws.Range("I113").Value = 1003
$ws.Range("J113").Value = 1293.6
$ws.Range("K113").Value = 1003
$ws.Range("L113").Value = 1293.6
$ws.Range("M113").Value = 1167
$ws.Range("N113").Value = -5633.6
$ws.Range("H126").Value = 2525.25
$ws.Range("I126").Value = 2431.077
$ws.Range("J126").Value = 2933.3333
$ws.Range("K126").Value = 7293.231000000001
$ws.Range("L126").Value = 8799.999899999999
$ws.Range("M126").Value = -4823.231000000001
$ws.Range("N126").Value = -13739.9999
$ws.Range("H132").Value = 11433.167
$ws.Range("I132").Value = 10044.333
$ws.Range("J132").Value = 15599.667
$ws.Range("K132").Value = 30132.999
$ws.Range("L132").Value = 46799.001
$ws.Range("M132").Value = -27602.999
$ws.Range("N132").Value = -51859.001
$ws.Range("H136").Value = 3214.9
$ws.Range("J136").Value = 3101.3333
$ws.Range("L136").Value = 9303.999899999999
$ws.Range("N136").Value = -14403.9999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 214.93333
$ws.Range("I11").Value = 87.85714
$ws.Range("K11").Value = 263.57142
$ws.Range("M11").Value = -123.57142

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37910.91
$ws.Range("J70").Value = 22815.666
$ws.Range("L70").Value = 22815.666
$ws.Range("N70").Value = -23355.666
$ws.Range("H73").Value = 37910.91
$ws.Range("J73").Value = 22815.666
$ws.Range("L73").Value = 22815.666
$ws.Range("N73").Value = -24687.666
$ws.Range("H80").Value = 2719
$ws.Range("I80").Value = 2648.75
$ws.Range("K80").Value = 2648.75
$ws.Range("M80").Value = -1650.75
$ws.Range("H83").Value = 2719
$ws.Range("I83").Value = 2648.75
$ws.Range("K83").Value = 13243.75
$ws.Range("M83").Value = -8251.75
$ws.Range("H102").Value = 4666.074
$ws.Range("I102").Value = 2299.6924
$ws.Range("K102").Value = 2299.6924
$ws.Range("M102").Value = -677.6923999999999
$ws.Range("H113").Value = 4314
$ws.Range("I113").Value = 4222.25
$ws.Range("J113").Value = 4497.5
$ws.Range("K113").Value = 4222.25
$ws.Range("L113").Value = 4497.5
$ws.Range("M113").Value = -2052.25
$ws.Range("N113").Value = -8837.5
$ws.Range("H122").Value = 3821.087
$ws.Range("I122").Value = 3803.3333
$ws.Range("J122").Value = 4007.5
$ws.Range("K122").Value = 11409.9999
$ws.Range("L122").Value = 12022.5
$ws.Range("M122").Value = -8959.999899999999
$ws.Range("N122").Value = -16922.5
$ws.Range("H132").Value = 4811.0156
$ws.Range("I132").Value = 4218.868
$ws.Range("K132").Value = 12656.604
$ws.Range("M132").Value = -10126.604

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15156357
$ws.Range("I7").Value = 4893
$ws.Range("J7").Value = 41671420
$ws.Range("K7").Value = 4893
$ws.Range("L7").Value = 41671420
$ws.Range("M7").Value = -4781
$ws.Range("N7").Value = -41671644
$ws.Range("H93").Value = 3440.6191
$ws.Range("I93").Value = 4122.3125
$ws.Range("J93").Value = 1259.2
$ws.Range("K93").Value = 4122.3125
$ws.Range("L93").Value = 1259.2
$ws.Range("M93").Value = -2874.3125
$ws.Range("N93").Value = -3755.2
$ws.Range("H126").Value = 15156357
$ws.Range("I126").Value = 4893
$ws.Range("J126").Value = 41671420
$ws.Range("K126").Value = 14679
$ws.Range("L126").Value = 125014260
$ws.Range("M126").Value = -12209
$ws.Range("N126").Value = -125019200
$ws.Range("H136").Value = 3115.6155
$ws.Range("I136").Value = 2124.125
$ws.Range("K136").Value = 6372.375
$ws.Range("M136").Value = -3822.375

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 22295
$ws.Range("J98").Value = 22295
$ws.Range("L98").Value = 22295
$ws.Range("N98").Value = -28285
$ws.Range("H107").Value = 862.5143
$ws.Range("J107").Value = 1350.3636
$ws.Range("L107").Value = 4051.0908
$ws.Range("N107").Value = -7891.0908
$ws.Range("H132").Value = 3044.889
$ws.Range("I132").Value = 1955.3549
$ws.Range("J132").Value = 9800
$ws.Range("K132").Value = 5866.0647
$ws.Range("L132").Value = 29400
$ws.Range("M132").Value = -3336.0647
$ws.Range("N132").Value = -34460
